$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_PAR"

$values = @(
    3.30551412936074,
    205.04199274182,
    180.872657085349,
    88.269466312663,
    0.337250969841363,
    207.014138318789,
    182.388227531804,
    48.1610250280693,
    0.0372677996249963,
    210.190123630547,
    179.678623671215,
    94.0334077516263,
    1.35000515117752,
    138.424386313528,
    172.193928226923,
    108.408842086369
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $values[$i]
}
